# Adds a new data row (row 33) to the master-reg_center_user_machine sheet,
# mirroring the existing rows' layout:
#   A: regcntr_id   B: usr_id   C: machine_id   D: lang_code
#   E: is_active    F: cr_by    G: cr_dtimes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 33

$ws.Cells.Item($newRow, 1).Value = 10002
$ws.Cells.Item($newRow, 2).Value = 110032
$ws.Cells.Item($newRow, 3).Value = 10032
$ws.Cells.Item($newRow, 4).Value = "eng"
$ws.Cells.Item($newRow, 5).Value = $true
$ws.Cells.Item($newRow, 6).Value = "superadmin"
$ws.Cells.Item($newRow, 7).Value = "now()"

# Keep the visible selection consistent with the recorded state after the edit.
$ws.Range("C30").Select() | Out-Null
